$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled name "H.BOWN" -> "H.BROWN" in column B (s2cDNAPreparer)
# for rows 2 through 19 (continuing the earlier fix already applied to column E).
$ws.Range("B2:B19").Value = "H.BROWN"

# Reflect the resulting active selection on column B.
$ws.Range("B2:B19").Select()
